# Fix typo in header cell D1: "Low_Compex" -> "Low_Complex"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").Value = "Low_Complex"
